$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Slide 6 table: switch the table style applied to the graphicFrame
#    table from the custom "Table_0" style to the built-in
#    {A82717E8-A00D-49E7-B7FE-838206CEBE75} style.
# ---------------------------------------------------------------------
$tableSlide = $p.Slides.Item(6)
for ($i = 1; $i -le $tableSlide.Shapes.Count; $i++) {
    $shp = $tableSlide.Shapes.Item($i)
    if ($shp.HasTable) {
        $shp.Table.ApplyStyle("{A82717E8-A00D-49E7-B7FE-838206CEBE75}")
    }
}

# ---------------------------------------------------------------------
# 2) Theme colors: the deck's active theme (the one used by the slide
#    master / all slides) switches from the "Integral" palette to the
#    stock "Office" palette. Re-point every theme color in the
#    12-entry scheme (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink)
#    to the corresponding Office color. Colors are VBA BGR integers.
# ---------------------------------------------------------------------
$officeColors = @(0, 16777215, 6968388, 15132391, 13998939, 3243501, 10855845, 49407, 12874308, 4697456, 12673797, 7491477)

$refSlide = $p.Slides.Item(1)
$colorScheme = $refSlide.ThemeColorScheme
for ($i = 1; $i -le $colorScheme.Count; $i++) {
    $colorScheme.Colors($i).RGB = $officeColors[$i - 1]
}
